# Generate Report for Handoff
#
# Updates the localization-status report to reflect a fresh handoff run:
#   - bumps the two timestamp values that are shared across the "Overview",
#     "zh-cn" and "de-de" sheets (the xliff-generate / handoff datetimes for
#     the batch of files that were just re-handed-off)
#   - marks those same rows' Priority column as "ht" (hot/high-priority
#     handoff) on both the "zh-cn" and "de-de" language sheets

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows (as they appear in each sheet) that belong to the batch of files
# that were just re-handed-off.
$rows = @(7, 8, 9, 10, 11, 13)

foreach ($r in $rows) {
    # Overview!G<r> : "Latest HO Xliff Generate Date" 2016-09-04 14:24:14 -> 14:24:33
    $overview.Cells.Item($r, 7).Value = "2016-09-04 14:24:33"

    # de-de!H<r> : "Latest Handoff Datetime" 2016-09-04 14:24:14 -> 14:24:33
    $dede.Cells.Item($r, 8).Value = "2016-09-04 14:24:33"

    # zh-cn!H<r> : "Latest Handoff Datetime" 2016-09-04 14:24:08 -> 14:24:29
    $zhcn.Cells.Item($r, 8).Value = "2016-09-04 14:24:29"

    # Priority column (E) flips from blank to "ht" on both language sheets
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 5).Value = "ht"
}
